$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The edit moves the "_GoBack" bookmark from the end of the Google
# Books URL paragraph (at the end of the "Natural language
# processing:" section) down to a brand-new final paragraph that
# holds a new SAGE reference URL, inserted right after "Arguments:".
# ------------------------------------------------------------------

# 1. Drop the old "_GoBack" bookmark wherever it currently sits.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Locate the "Arguments:" paragraph; the paragraph immediately
#    after it is the trailing empty paragraph that should become the
#    new "https://methods-sagepub-com..." paragraph.
$targetPara = $null
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd("`r", "`n", "`a") -eq "Arguments:") {
        $targetPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($targetPara -eq $null) {
    # Fallback: just use the very last paragraph of the document.
    $targetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
}

# 3. Re-create "_GoBack" on that (still empty) paragraph first, then
#    insert the new URL text *before* the bookmark's position so the
#    bookmark ends up wrapping/following the newly typed run, exactly
#    like the original "<run/><bookmarkStart/><bookmarkEnd/>" shape.
$d.Bookmarks.Add("_GoBack", $targetPara.Range)
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Range.InsertBefore("https://methods-sagepub-com.ezproxy.napier.ac.uk/Reference/the-sage-encyclopedia-of-communication-research-methods/i1874.xml")

Write-Output "Paragraphs: $($d.Paragraphs.Count); GoBack exists: $($d.Bookmarks.Exists('_GoBack'))"
